$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1971722498084585
$ws.Range("C2").Value = 0.001640823957704524
$ws.Range("B3").Value = 0.3870545671486926
$ws.Range("C3").Value = 0.002521359895798399
